# ISD lab 2b fix
# Applies the row/style corrections described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content updates -------------------------------------------------
# E2: "теория" -> "Задачи на листе взять"
$ws.Range("E2").Value = "Задачи на листе взять"

# E4 used to hold "default, constraint, pattern переписать 2а" - clear it
$ws.Range("E4").Value = $null

# C5 used to hold "отчёты" - clear it
$ws.Range("C5").Value = $null

# --- Fill colour updates ----------------------------------------------------
# Excel's Interior.Color / Interior.PatternColor use BGR-encoded integers.
# Reference colours (ARGB hex -> BGR decimal):
#   FF81D41A / FF92D050  (reversed green)   -> 1758337 / 5296274
#   FFFF0000 / FF993300  (red / brown)      -> 255     / 13209
#   FF000000 / FF003300  (black / dk green) -> 0       / 13056
#   FFFFFF00 / FFFFFF00  (yellow)           -> 65535   / 65535
#   FFFFBF00 / FFFF9900  (orange)           -> 49151   / 39423

function Set-Fill($addr, $fg, $bg) {
    $c = $ws.Range($addr)
    $c.Interior.Pattern = 1
    $c.Interior.Color = $fg
    $c.Interior.PatternColor = $bg
}

# Cells that become the reversed-green fill
foreach ($addr in @("E2", "C5", "D5", "F8", "G8")) {
    Set-Fill $addr 1758337 5296274
}

# Cells that become red/brown
foreach ($addr in @("D3", "E3")) {
    Set-Fill $addr 255 13209
}

# Cell that becomes the new black / dark-green fill
Set-Fill "I4" 0 13056

# Cell that becomes yellow
Set-Fill "D7" 65535 65535

# Cell that becomes orange
Set-Fill "E7" 49151 39423

# --- View state (best effort) ----------------------------------------------
$ws.Activate()
try { $excel.ActiveWindow.ScrollRow = 4 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
$ws.Range("D5").Select()

Write-Host "ISD lab 2b fix applied"
